{"js": "// Replace the two-digit-divided-by-one-digit problem texts throughout\n// the document body (including inside the table cells) with their\n// updated values, as described by the commit diff.\nconst replacements = [\n  [\"91\u00f72=\", \"18\u00f73=\"],\n  [\"52\u00f75=\", \"67\u00f72=\"],\n  [\"89\u00f77=\", \"61\u00f72=\"],\n  [\"60\u00f72=\", \"44\u00f78=\"],\n  [\"34\u00f75=\", \"10\u00f77=\"],\n  [\"77\u00f77=\", \"23\u00f74=\"],\n  [\"68\u00f74=\", \"42\u00f73=\"],\n  [\"71\u00f77=\", \"26\u00f72=\"],\n  [\"30\u00f79=\", \"60\u00f75=\"],\n  [\"96\u00f72=\", \"78\u00f73=\"],\n  [\"50\u00f74=\", \"71\u00f76=\"],\n  [\"97\u00f76=\", \"39\u00f76=\"],\n  [\"51\u00f76=\", \"25\u00f72=\"],\n  [\"58\u00f73=\", \"17\u00f72=\"],\n  [\"32\u00f78=\", \"31\u00f74=\"],\n  [\"25\u00f75=\", \"33\u00f72=\"],\n  [\"24\u00f72=\", \"24\u00f73=\"],\n  [\"31\u00f77=\", \"72\u00f79=\"],\n  [\"71\u00f79=\", \"44\u00f79=\"],\n  [\"11\u00f78=\", \"49\u00f77=\"],\n  [\"22\u00f78=\", \"59\u00f78=\"],\n  [\"28\u00f72=\", \"71\u00f72=\"],\n  [\"94\u00f72=\", \"55\u00f78=\"],\n  [\"21\u00f78=\", \"15\u00f77=\"],\n  [\"62\u00f78=\", \"57\u00f77=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const found = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the two-digit-divided-by-one-digit problem texts throughout\n# the document body (including inside the table cells) with their\n# updated values, as described by the commit diff.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"91\u00f72=\", \"18\u00f73=\"),\n    @(\"52\u00f75=\", \"67\u00f72=\"),\n    @(\"89\u00f77=\", \"61\u00f72=\"),\n    @(\"60\u00f72=\", \"44\u00f78=\"),\n    @(\"34\u00f75=\", \"10\u00f77=\"),\n    @(\"77\u00f77=\", \"23\u00f74=\"),\n    @(\"68\u00f74=\", \"42\u00f73=\"),\n    @(\"71\u00f77=\", \"26\u00f72=\"),\n    @(\"30\u00f79=\", \"60\u00f75=\"),\n    @(\"96\u00f72=\", \"78\u00f73=\"),\n    @(\"50\u00f74=\", \"71\u00f76=\"),\n    @(\"97\u00f76=\", \"39\u00f76=\"),\n    @(\"51\u00f76=\", \"25\u00f72=\"),\n    @(\"58\u00f73=\", \"17\u00f72=\"),\n    @(\"32\u00f78=\", \"31\u00f74=\"),\n    @(\"25\u00f75=\", \"33\u00f72=\"),\n    @(\"24\u00f72=\", \"24\u00f73=\"),\n    @(\"31\u00f77=\", \"72\u00f79=\"),\n    @(\"71\u00f79=\", \"44\u00f79=\"),\n    @(\"11\u00f78=\", \"49\u00f77=\"),\n    @(\"22\u00f78=\", \"59\u00f78=\"),\n    @(\"28\u00f72=\", \"71\u00f72=\"),\n    @(\"94\u00f72=\", \"55\u00f78=\"),\n    @(\"21\u00f78=\", \"15\u00f77=\"),\n    @(\"62\u00f78=\", \"57\u00f77=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute(\n        $oldText,   # FindText\n        $true,      # MatchCase\n        $false,     # MatchWholeWord\n        $false,     # MatchWildcards\n        $false,     # MatchSoundsLike\n        $false,     # MatchAllWordForms\n        $true,      # Forward\n        1,          # Wrap (wdFindContinue)\n        $false,     # Format\n        $newText,   # ReplaceWith\n        2           # Replace (wdReplaceAll)\n    )\n}\n"}
